$wb = $excel.ActiveWorkbook

# ---- Triple Dribble ----
$ws = $wb.Worksheets.Item("Triple Dribble")
# row 53
$ws.Range("A4:N4").Copy()
$ws.Range("A53:N53").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G53").PasteSpecial(-4122)
$ws.Cells.Item(53,1).Value = "SHADE"
$ws.Cells.Item(53,2).Value = "MOE"
$ws.Cells.Item(53,3).Value = "WILLOW"
$ws.Cells.Item(53,4).Value = "MEEPLE"
$ws.Cells.Item(53,5).Value = "LARRY & LAWRIE"
$ws.Cells.Item(53,6).Value = "KAZE"
$ws.Cells.Item(53,7).Value = "Equipo 1"
$ws.Cells.Item(53,8).Value = "MM"
$ws.Cells.Item(53,9).Value = "RC|Battoman"
$ws.Cells.Item(53,10).Value = "RC|Shu"
$ws.Cells.Item(53,11).Value = "CR|Moya"
$ws.Cells.Item(53,12).Value = "CR|Milkreo"
$ws.Cells.Item(53,13).Value = "CR|Tensai"
$ws.Cells.Item(53,14).Value = "20250725T135359.000Z"

# row 54
$ws.Range("A4:N4").Copy()
$ws.Range("A54:N54").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G54").PasteSpecial(-4122)
$ws.Cells.Item(54,1).Value = "SHADE"
$ws.Cells.Item(54,2).Value = "MOE"
$ws.Cells.Item(54,3).Value = "WILLOW"
$ws.Cells.Item(54,4).Value = "MEEPLE"
$ws.Cells.Item(54,5).Value = "LARRY & LAWRIE"
$ws.Cells.Item(54,6).Value = "KAZE"
$ws.Cells.Item(54,7).Value = "Equipo 2"
$ws.Cells.Item(54,8).Value = "MM"
$ws.Cells.Item(54,9).Value = "RC|Battoman"
$ws.Cells.Item(54,10).Value = "RC|Shu"
$ws.Cells.Item(54,11).Value = "CR|Moya"
$ws.Cells.Item(54,12).Value = "CR|Milkreo"
$ws.Cells.Item(54,13).Value = "CR|Tensai"
$ws.Cells.Item(54,14).Value = "20250725T135055.000Z"

# row 55
$ws.Range("A4:N4").Copy()
$ws.Range("A55:N55").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G55").PasteSpecial(-4122)
$ws.Cells.Item(55,1).Value = "SHADE"
$ws.Cells.Item(55,2).Value = "MOE"
$ws.Cells.Item(55,3).Value = "WILLOW"
$ws.Cells.Item(55,4).Value = "MEEPLE"
$ws.Cells.Item(55,5).Value = "LARRY & LAWRIE"
$ws.Cells.Item(55,6).Value = "KAZE"
$ws.Cells.Item(55,7).Value = "Equipo 1"
$ws.Cells.Item(55,8).Value = "MM"
$ws.Cells.Item(55,9).Value = "RC|Battoman"
$ws.Cells.Item(55,10).Value = "RC|Shu"
$ws.Cells.Item(55,11).Value = "CR|Moya"
$ws.Cells.Item(55,12).Value = "CR|Milkreo"
$ws.Cells.Item(55,13).Value = "CR|Tensai"
$ws.Cells.Item(55,14).Value = "20250725T134922.000Z"


# ---- Double Swoosh ----
$ws = $wb.Worksheets.Item("Double Swoosh")
# row 14
$ws.Range("A4:N4").Copy()
$ws.Range("A14:N14").PasteSpecial(-4122)
$ws.Range("G11").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Cells.Item(14,1).Value = "GRAY"
$ws.Cells.Item(14,2).Value = "LILY"
$ws.Cells.Item(14,3).Value = "EMZ"
$ws.Cells.Item(14,4).Value = "KAZE"
$ws.Cells.Item(14,5).Value = "TARA"
$ws.Cells.Item(14,6).Value = "SANDY"
$ws.Cells.Item(14,7).Value = "Equipo 2"
$ws.Cells.Item(14,8).Value = "GEN|Moding"
$ws.Cells.Item(14,9).Value = "GEN|BONOX2"
$ws.Cells.Item(14,10).Value = "GEN|cookie"
$ws.Cells.Item(14,11).Value = "FZ|Toridesu"
$ws.Cells.Item(14,12).Value = "FZ|Danshari"
$ws.Cells.Item(14,13).Value = "FZ|Mira"
$ws.Cells.Item(14,14).Value = "20250725T133435.000Z"

# row 15
$ws.Range("A4:N4").Copy()
$ws.Range("A15:N15").PasteSpecial(-4122)
$ws.Range("G11").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Cells.Item(15,1).Value = "GRAY"
$ws.Cells.Item(15,2).Value = "LILY"
$ws.Cells.Item(15,3).Value = "EMZ"
$ws.Cells.Item(15,4).Value = "KAZE"
$ws.Cells.Item(15,5).Value = "TARA"
$ws.Cells.Item(15,6).Value = "SANDY"
$ws.Cells.Item(15,7).Value = "Equipo 2"
$ws.Cells.Item(15,8).Value = "GEN|Moding"
$ws.Cells.Item(15,9).Value = "GEN|BONOX2"
$ws.Cells.Item(15,10).Value = "GEN|cookie"
$ws.Cells.Item(15,11).Value = "FZ|Toridesu"
$ws.Cells.Item(15,12).Value = "FZ|Danshari"
$ws.Cells.Item(15,13).Value = "FZ|Mira"
$ws.Cells.Item(15,14).Value = "20250725T133236.000Z"

# row 16
$ws.Range("A4:N4").Copy()
$ws.Range("A16:N16").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value = "GRAY"
$ws.Cells.Item(16,2).Value = "LILY"
$ws.Cells.Item(16,3).Value = "EMZ"
$ws.Cells.Item(16,4).Value = "KAZE"
$ws.Cells.Item(16,5).Value = "TARA"
$ws.Cells.Item(16,6).Value = "SANDY"
$ws.Cells.Item(16,7).Value = "Equipo 1"
$ws.Cells.Item(16,8).Value = "GEN|Moding"
$ws.Cells.Item(16,9).Value = "GEN|BONOX2"
$ws.Cells.Item(16,10).Value = "GEN|cookie"
$ws.Cells.Item(16,11).Value = "FZ|Toridesu"
$ws.Cells.Item(16,12).Value = "FZ|Danshari"
$ws.Cells.Item(16,13).Value = "FZ|Mira"
$ws.Cells.Item(16,14).Value = "20250725T133038.000Z"


# ---- Crystal Arcade ----
$ws = $wb.Worksheets.Item("Crystal Arcade")
# row 65
$ws.Range("A4:N4").Copy()
$ws.Range("A65:N65").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G65").PasteSpecial(-4122)
$ws.Cells.Item(65,1).Value = "KAZE"
$ws.Cells.Item(65,2).Value = "LOU"
$ws.Cells.Item(65,3).Value = "CROW"
$ws.Cells.Item(65,4).Value = "KIT"
$ws.Cells.Item(65,5).Value = "DRACO"
$ws.Cells.Item(65,6).Value = "JANET"
$ws.Cells.Item(65,7).Value = "Equipo 1"
$ws.Cells.Item(65,8).Value = "NAVI|Ryohei"
$ws.Cells.Item(65,9).Value = "NAVI|Achapi"
$ws.Cells.Item(65,10).Value = "NAVI|Kuru"
$ws.Cells.Item(65,11).Value = "あの頃のしてたんぽ👍"
$ws.Cells.Item(65,12).Value = "あの頃のしずく👍"
$ws.Cells.Item(65,13).Value = "ZETA|Levi"
$ws.Cells.Item(65,14).Value = "20250725T135312.000Z"

# row 66
$ws.Range("A4:N4").Copy()
$ws.Range("A66:N66").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G66").PasteSpecial(-4122)
$ws.Cells.Item(66,1).Value = "KAZE"
$ws.Cells.Item(66,2).Value = "LOU"
$ws.Cells.Item(66,3).Value = "CROW"
$ws.Cells.Item(66,4).Value = "KIT"
$ws.Cells.Item(66,5).Value = "DRACO"
$ws.Cells.Item(66,6).Value = "JANET"
$ws.Cells.Item(66,7).Value = "Equipo 1"
$ws.Cells.Item(66,8).Value = "NAVI|Ryohei"
$ws.Cells.Item(66,9).Value = "NAVI|Achapi"
$ws.Cells.Item(66,10).Value = "NAVI|Kuru"
$ws.Cells.Item(66,11).Value = "あの頃のしてたんぽ👍"
$ws.Cells.Item(66,12).Value = "あの頃のしずく👍"
$ws.Cells.Item(66,13).Value = "ZETA|Levi"
$ws.Cells.Item(66,14).Value = "20250725T134954.000Z"

# row 67
$ws.Range("A4:N4").Copy()
$ws.Range("A67:N67").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G67").PasteSpecial(-4122)
$ws.Cells.Item(67,1).Value = "HANK"
$ws.Cells.Item(67,2).Value = "SPIKE"
$ws.Cells.Item(67,3).Value = "JAE-YONG"
$ws.Cells.Item(67,4).Value = "EMZ"
$ws.Cells.Item(67,5).Value = "KIT"
$ws.Cells.Item(67,6).Value = "GUS"
$ws.Cells.Item(67,7).Value = "Equipo 2"
$ws.Cells.Item(67,8).Value = "NAVI|Ryohei"
$ws.Cells.Item(67,9).Value = "NAVI|Achapi"
$ws.Cells.Item(67,10).Value = "NAVI|Kuru"
$ws.Cells.Item(67,11).Value = "あの頃のしてたんぽ👍"
$ws.Cells.Item(67,12).Value = "あの頃のしずく👍"
$ws.Cells.Item(67,13).Value = "ZETA|Levi"
$ws.Cells.Item(67,14).Value = "20250725T134347.000Z"

# row 68
$ws.Range("A4:N4").Copy()
$ws.Range("A68:N68").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G68").PasteSpecial(-4122)
$ws.Cells.Item(68,1).Value = "HANK"
$ws.Cells.Item(68,2).Value = "SPIKE"
$ws.Cells.Item(68,3).Value = "JAE-YONG"
$ws.Cells.Item(68,4).Value = "EMZ"
$ws.Cells.Item(68,5).Value = "KIT"
$ws.Cells.Item(68,6).Value = "GUS"
$ws.Cells.Item(68,7).Value = "Equipo 2"
$ws.Cells.Item(68,8).Value = "NAVI|Ryohei"
$ws.Cells.Item(68,9).Value = "NAVI|Achapi"
$ws.Cells.Item(68,10).Value = "NAVI|Kuru"
$ws.Cells.Item(68,11).Value = "あの頃のしてたんぽ👍"
$ws.Cells.Item(68,12).Value = "あの頃のしずく👍"
$ws.Cells.Item(68,13).Value = "ZETA|Levi"
$ws.Cells.Item(68,14).Value = "20250725T134103.000Z"

# row 69
$ws.Range("A4:N4").Copy()
$ws.Range("A69:N69").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G69").PasteSpecial(-4122)
$ws.Cells.Item(69,1).Value = "SURGE"
$ws.Cells.Item(69,2).Value = "SPIKE"
$ws.Cells.Item(69,3).Value = "LILY"
$ws.Cells.Item(69,4).Value = "CORDELIUS"
$ws.Cells.Item(69,5).Value = "BULL"
$ws.Cells.Item(69,6).Value = "CHARLIE"
$ws.Cells.Item(69,7).Value = "Equipo 2"
$ws.Cells.Item(69,8).Value = "MM"
$ws.Cells.Item(69,9).Value = "RC|Battoman"
$ws.Cells.Item(69,10).Value = "RC|Shu"
$ws.Cells.Item(69,11).Value = "CR|Milkreo"
$ws.Cells.Item(69,12).Value = "CR|Moya"
$ws.Cells.Item(69,13).Value = "CR|Tensai"
$ws.Cells.Item(69,14).Value = "20250725T134230.000Z"

# row 70
$ws.Range("A4:N4").Copy()
$ws.Range("A70:N70").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G70").PasteSpecial(-4122)
$ws.Cells.Item(70,1).Value = "SURGE"
$ws.Cells.Item(70,2).Value = "SPIKE"
$ws.Cells.Item(70,3).Value = "LILY"
$ws.Cells.Item(70,4).Value = "CORDELIUS"
$ws.Cells.Item(70,5).Value = "BULL"
$ws.Cells.Item(70,6).Value = "CHARLIE"
$ws.Cells.Item(70,7).Value = "Equipo 2"
$ws.Cells.Item(70,8).Value = "MM"
$ws.Cells.Item(70,9).Value = "RC|Battoman"
$ws.Cells.Item(70,10).Value = "RC|Shu"
$ws.Cells.Item(70,11).Value = "CR|Milkreo"
$ws.Cells.Item(70,12).Value = "CR|Moya"
$ws.Cells.Item(70,13).Value = "CR|Tensai"
$ws.Cells.Item(70,14).Value = "20250725T134014.000Z"

# row 71
$ws.Range("A4:N4").Copy()
$ws.Range("A71:N71").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G71").PasteSpecial(-4122)
$ws.Cells.Item(71,1).Value = "CORDELIUS"
$ws.Cells.Item(71,2).Value = "DARRYL"
$ws.Cells.Item(71,3).Value = "JANET"
$ws.Cells.Item(71,4).Value = "BUSTER"
$ws.Cells.Item(71,5).Value = "GUS"
$ws.Cells.Item(71,6).Value = "LUMI"
$ws.Cells.Item(71,7).Value = "Equipo 2"
$ws.Cells.Item(71,8).Value = "MM"
$ws.Cells.Item(71,9).Value = "RC|Shu"
$ws.Cells.Item(71,10).Value = "RC|Battoman"
$ws.Cells.Item(71,11).Value = "CR|Moya"
$ws.Cells.Item(71,12).Value = "CR|Milkreo"
$ws.Cells.Item(71,13).Value = "CR|Tensai"
$ws.Cells.Item(71,14).Value = "20250725T133349.000Z"

# row 72
$ws.Range("A4:N4").Copy()
$ws.Range("A72:N72").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Range("G72").PasteSpecial(-4122)
$ws.Cells.Item(72,1).Value = "CORDELIUS"
$ws.Cells.Item(72,2).Value = "DARRYL"
$ws.Cells.Item(72,3).Value = "JANET"
$ws.Cells.Item(72,4).Value = "BUSTER"
$ws.Cells.Item(72,5).Value = "GUS"
$ws.Cells.Item(72,6).Value = "LUMI"
$ws.Cells.Item(72,7).Value = "Equipo 2"
$ws.Cells.Item(72,8).Value = "MM"
$ws.Cells.Item(72,9).Value = "RC|Shu"
$ws.Cells.Item(72,10).Value = "RC|Battoman"
$ws.Cells.Item(72,11).Value = "CR|Moya"
$ws.Cells.Item(72,12).Value = "CR|Milkreo"
$ws.Cells.Item(72,13).Value = "CR|Tensai"
$ws.Cells.Item(72,14).Value = "20250725T133034.000Z"


# ---- New Horizons ----
$ws = $wb.Worksheets.Item("New Horizons")
# row 78
$ws.Range("A4:N4").Copy()
$ws.Range("A78:N78").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G78").PasteSpecial(-4122)
$ws.Cells.Item(78,1).Value = "MR. P"
$ws.Cells.Item(78,2).Value = "CHARLIE"
$ws.Cells.Item(78,3).Value = "HANK"
$ws.Cells.Item(78,4).Value = "MEEPLE"
$ws.Cells.Item(78,5).Value = "BELLE"
$ws.Cells.Item(78,6).Value = "GRAY"
$ws.Cells.Item(78,7).Value = "Equipo 2"
$ws.Cells.Item(78,8).Value = "NAVI|Achapi"
$ws.Cells.Item(78,9).Value = "NAVI|Kuru"
$ws.Cells.Item(78,10).Value = "NAVI|Ryohei"
$ws.Cells.Item(78,11).Value = "あの頃のしてたんぽ👍"
$ws.Cells.Item(78,12).Value = "ZETA|Levi"
$ws.Cells.Item(78,13).Value = "あの頃のしずく👍"
$ws.Cells.Item(78,14).Value = "20250725T133446.000Z"

# row 79
$ws.Range("A4:N4").Copy()
$ws.Range("A79:N79").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G79").PasteSpecial(-4122)
$ws.Cells.Item(79,1).Value = "MR. P"
$ws.Cells.Item(79,2).Value = "CHARLIE"
$ws.Cells.Item(79,3).Value = "HANK"
$ws.Cells.Item(79,4).Value = "MEEPLE"
$ws.Cells.Item(79,5).Value = "BELLE"
$ws.Cells.Item(79,6).Value = "GRAY"
$ws.Cells.Item(79,7).Value = "Equipo 1"
$ws.Cells.Item(79,8).Value = "NAVI|Achapi"
$ws.Cells.Item(79,9).Value = "NAVI|Kuru"
$ws.Cells.Item(79,10).Value = "NAVI|Ryohei"
$ws.Cells.Item(79,11).Value = "あの頃のしてたんぽ👍"
$ws.Cells.Item(79,12).Value = "ZETA|Levi"
$ws.Cells.Item(79,13).Value = "あの頃のしずく👍"
$ws.Cells.Item(79,14).Value = "20250725T133209.000Z"

# row 80
$ws.Range("A4:N4").Copy()
$ws.Range("A80:N80").PasteSpecial(-4122)
$ws.Range("G8").Copy()
$ws.Range("G80").PasteSpecial(-4122)
$ws.Cells.Item(80,1).Value = "MR. P"
$ws.Cells.Item(80,2).Value = "CHARLIE"
$ws.Cells.Item(80,3).Value = "HANK"
$ws.Cells.Item(80,4).Value = "MEEPLE"
$ws.Cells.Item(80,5).Value = "BELLE"
$ws.Cells.Item(80,6).Value = "GRAY"
$ws.Cells.Item(80,7).Value = "Equipo 2"
$ws.Cells.Item(80,8).Value = "NAVI|Achapi"
$ws.Cells.Item(80,9).Value = "NAVI|Kuru"
$ws.Cells.Item(80,10).Value = "NAVI|Ryohei"
$ws.Cells.Item(80,11).Value = "あの頃のしてたんぽ👍"
$ws.Cells.Item(80,12).Value = "ZETA|Levi"
$ws.Cells.Item(80,13).Value = "あの頃のしずく👍"
$ws.Cells.Item(80,14).Value = "20250725T133026.000Z"


# ---- Ring of Fire ----
$ws = $wb.Worksheets.Item("Ring of Fire")
# row 45
$ws.Range("A4:N4").Copy()
$ws.Range("A45:N45").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G45").PasteSpecial(-4122)
$ws.Cells.Item(45,1).Value = "LOLA"
$ws.Cells.Item(45,2).Value = "AMBER"
$ws.Cells.Item(45,3).Value = "GRIFF"
$ws.Cells.Item(45,4).Value = "MELODIE"
$ws.Cells.Item(45,5).Value = "DRACO"
$ws.Cells.Item(45,6).Value = "PAM"
$ws.Cells.Item(45,7).Value = "Equipo 1"
$ws.Cells.Item(45,8).Value = "GEN|cookie"
$ws.Cells.Item(45,9).Value = "GEN|BONOX2"
$ws.Cells.Item(45,10).Value = "GEN|Moding"
$ws.Cells.Item(45,11).Value = "FZ|Toridesu"
$ws.Cells.Item(45,12).Value = "FZ|Danshari"
$ws.Cells.Item(45,13).Value = "FZ|Mira"
$ws.Cells.Item(45,14).Value = "20250725T135440.000Z"

# row 46
$ws.Range("A4:N4").Copy()
$ws.Range("A46:N46").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G46").PasteSpecial(-4122)
$ws.Cells.Item(46,1).Value = "BEA"
$ws.Cells.Item(46,2).Value = "PAM"
$ws.Cells.Item(46,3).Value = "FRANK"
$ws.Cells.Item(46,4).Value = "DRACO"
$ws.Cells.Item(46,5).Value = "PENNY"
$ws.Cells.Item(46,6).Value = "BO"
$ws.Cells.Item(46,7).Value = "Equipo 2"
$ws.Cells.Item(46,8).Value = "GEN|cookie"
$ws.Cells.Item(46,9).Value = "GEN|Moding"
$ws.Cells.Item(46,10).Value = "GEN|BONOX2"
$ws.Cells.Item(46,11).Value = "FZ|Toridesu"
$ws.Cells.Item(46,12).Value = "FZ|Danshari"
$ws.Cells.Item(46,13).Value = "FZ|Mira"
$ws.Cells.Item(46,14).Value = "20250725T134642.000Z"

# row 47
$ws.Range("A4:N4").Copy()
$ws.Range("A47:N47").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G47").PasteSpecial(-4122)
$ws.Cells.Item(47,1).Value = "BEA"
$ws.Cells.Item(47,2).Value = "PAM"
$ws.Cells.Item(47,3).Value = "FRANK"
$ws.Cells.Item(47,4).Value = "DRACO"
$ws.Cells.Item(47,5).Value = "PENNY"
$ws.Cells.Item(47,6).Value = "BO"
$ws.Cells.Item(47,7).Value = "Equipo 1"
$ws.Cells.Item(47,8).Value = "GEN|cookie"
$ws.Cells.Item(47,9).Value = "GEN|Moding"
$ws.Cells.Item(47,10).Value = "GEN|BONOX2"
$ws.Cells.Item(47,11).Value = "FZ|Toridesu"
$ws.Cells.Item(47,12).Value = "FZ|Danshari"
$ws.Cells.Item(47,13).Value = "FZ|Mira"
$ws.Cells.Item(47,14).Value = "20250725T134342.000Z"

# row 48
$ws.Range("A4:N4").Copy()
$ws.Range("A48:N48").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G48").PasteSpecial(-4122)
$ws.Cells.Item(48,1).Value = "BEA"
$ws.Cells.Item(48,2).Value = "PAM"
$ws.Cells.Item(48,3).Value = "FRANK"
$ws.Cells.Item(48,4).Value = "DRACO"
$ws.Cells.Item(48,5).Value = "PENNY"
$ws.Cells.Item(48,6).Value = "BO"
$ws.Cells.Item(48,7).Value = "Equipo 2"
$ws.Cells.Item(48,8).Value = "GEN|cookie"
$ws.Cells.Item(48,9).Value = "GEN|Moding"
$ws.Cells.Item(48,10).Value = "GEN|BONOX2"
$ws.Cells.Item(48,11).Value = "FZ|Toridesu"
$ws.Cells.Item(48,12).Value = "FZ|Danshari"
$ws.Cells.Item(48,13).Value = "FZ|Mira"
$ws.Cells.Item(48,14).Value = "20250725T134134.000Z"


$excel.CutCopyMode = 0
Write-Host "Done updating scrims_actualizado.xlsx"